$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the "location" header to "name" and the "network" header to "netid"
$ws.Range("B2").Value = "name"
$ws.Range("C2").Value = "netid"

# Move the active selection to B3
$ws.Range("B3").Select()
